$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Dhh"
$ws.Cells.Item(2,3).Value = "Boc"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 1.853892
$ws.Cells.Item(2,8).Value = 5.561676
$ws.Cells.Item(2,9).Value = 0.5711238486747862
$ws.Cells.Item(2,10).Value = 0.571123848674786
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 2.173625666666667
$ws.Cells.Item(2,14).Value = 6.520877
$ws.Cells.Item(2,15).Value = 0.03242430220886856
$ws.Cells.Item(2,16).Value = 0.03242430220886856
$ws.Cells.Item(2,17).Value = 4.029667234428
$ws.Cells.Item(2,18).Value = 36.26700510985201
$ws.Cells.Item(2,19).Value = 0.01851829226812338
$ws.Cells.Item(2,20).Value = 0.01851829226812338

$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Dhh"
$ws.Cells.Item(3,3).Value = "Boc"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 1.853892
$ws.Cells.Item(3,8).Value = 5.561676
$ws.Cells.Item(3,9).Value = 0.5711238486747862
$ws.Cells.Item(3,10).Value = 0.571123848674786
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 59.528614
$ws.Cells.Item(3,14).Value = 178.585842
$ws.Cells.Item(3,15).Value = 0.8879973217150474
$ws.Cells.Item(3,16).Value = 0.8879973217150473
$ws.Cells.Item(3,17).Value = 110.359621265688
$ws.Cells.Item(3,18).Value = 993.2365913911922
$ws.Cells.Item(3,19).Value = 0.5071564479908002
$ws.Cells.Item(3,20).Value = 0.5071564479908

$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Dhh"
$ws.Cells.Item(4,3).Value = "Boc"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 1.853892
$ws.Cells.Item(4,8).Value = 5.561676
$ws.Cells.Item(4,9).Value = 0.5711238486747862
$ws.Cells.Item(4,10).Value = 0.571123848674786
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 5.334689999999999
$ws.Cells.Item(4,14).Value = 16.00407
$ws.Cells.Item(4,15).Value = 0.0795783760760841
$ws.Cells.Item(4,16).Value = 0.0795783760760841
$ws.Cells.Item(4,17).Value = 9.889939113479999
$ws.Cells.Item(4,18).Value = 89.00945202132
$ws.Cells.Item(4,19).Value = 0.04544910841586268
$ws.Cells.Item(4,20).Value = 0.04544910841586267

$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Dhh"
$ws.Cells.Item(5,3).Value = "Boc"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.4262446666666667
$ws.Cells.Item(5,8).Value = 1.278734
$ws.Cells.Item(5,9).Value = 0.1313121230922664
$ws.Cells.Item(5,10).Value = 0.1313121230922664
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 2.173625666666667
$ws.Cells.Item(5,14).Value = 6.520877
$ws.Cells.Item(5,15).Value = 0.03242430220886856
$ws.Cells.Item(5,16).Value = 0.03242430220886856
$ws.Cells.Item(5,17).Value = 0.9264963477464444
$ws.Cells.Item(5,18).Value = 8.338467129718001
$ws.Cells.Item(5,19).Value = 0.004257703962831795
$ws.Cells.Item(5,20).Value = 0.004257703962831794

$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Dhh"
$ws.Cells.Item(6,3).Value = "Boc"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.4262446666666667
$ws.Cells.Item(6,8).Value = 1.278734
$ws.Cells.Item(6,9).Value = 0.1313121230922664
$ws.Cells.Item(6,10).Value = 0.1313121230922664
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 59.528614
$ws.Cells.Item(6,14).Value = 178.585842
$ws.Cells.Item(6,15).Value = 0.8879973217150474
$ws.Cells.Item(6,16).Value = 0.8879973217150473
$ws.Cells.Item(6,17).Value = 25.37375423155867
$ws.Cells.Item(6,18).Value = 228.363788084028
$ws.Cells.Item(6,19).Value = 0.1166048136146492
$ws.Cells.Item(6,20).Value = 0.1166048136146492

$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Dhh"
$ws.Cells.Item(7,3).Value = "Boc"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.4262446666666667
$ws.Cells.Item(7,8).Value = 1.278734
$ws.Cells.Item(7,9).Value = 0.1313121230922664
$ws.Cells.Item(7,10).Value = 0.1313121230922664
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 5.334689999999999
$ws.Cells.Item(7,14).Value = 16.00407
$ws.Cells.Item(7,15).Value = 0.0795783760760841
$ws.Cells.Item(7,16).Value = 0.0795783760760841
$ws.Cells.Item(7,17).Value = 2.27388316082
$ws.Cells.Item(7,18).Value = 20.46494844738
$ws.Cells.Item(7,19).Value = 0.01044960551478543
$ws.Cells.Item(7,20).Value = 0.01044960551478542

$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Dhh"
$ws.Cells.Item(8,3).Value = "Boc"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.9659053333333333
$ws.Cells.Item(8,8).Value = 2.897716
$ws.Cells.Item(8,9).Value = 0.2975640282329475
$ws.Cells.Item(8,10).Value = 0.2975640282329475
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 2.173625666666667
$ws.Cells.Item(8,14).Value = 6.520877
$ws.Cells.Item(8,15).Value = 0.03242430220886856
$ws.Cells.Item(8,16).Value = 0.03242430220886856
$ws.Cells.Item(8,17).Value = 2.099516624103555
$ws.Cells.Item(8,18).Value = 18.895649616932
$ws.Cells.Item(8,19).Value = 0.009648305977913387
$ws.Cells.Item(8,20).Value = 0.009648305977913387

$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Dhh"
$ws.Cells.Item(9,3).Value = "Boc"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.9659053333333333
$ws.Cells.Item(9,8).Value = 2.897716
$ws.Cells.Item(9,9).Value = 0.2975640282329475
$ws.Cells.Item(9,10).Value = 0.2975640282329475
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 59.528614
$ws.Cells.Item(9,14).Value = 178.585842
$ws.Cells.Item(9,15).Value = 0.8879973217150474
$ws.Cells.Item(9,16).Value = 0.8879973217150473
$ws.Cells.Item(9,17).Value = 57.49900574854134
$ws.Cells.Item(9,18).Value = 517.491051736872
$ws.Cells.Item(9,19).Value = 0.2642360601095981
$ws.Cells.Item(9,20).Value = 0.2642360601095981

$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Dhh"
$ws.Cells.Item(10,3).Value = "Boc"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.9659053333333333
$ws.Cells.Item(10,8).Value = 2.897716
$ws.Cells.Item(10,9).Value = 0.2975640282329475
$ws.Cells.Item(10,10).Value = 0.2975640282329475
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 5.334689999999999
$ws.Cells.Item(10,14).Value = 16.00407
$ws.Cells.Item(10,15).Value = 0.0795783760760841
$ws.Cells.Item(10,16).Value = 0.0795783760760841
$ws.Cells.Item(10,17).Value = 5.152805522679999
$ws.Cells.Item(10,18).Value = 46.37524970411999
$ws.Cells.Item(10,19).Value = 0.023679662145436
$ws.Cells.Item(10,20).Value = 0.023679662145436

